# The M2Doc field `{ m:'doc.html'.fromHTMLURI() }` is stored as a real Word
# field (fldChar begin / instrText tokens / fldChar end). The parser now
# expects that same field rewritten as plain literal text (one run per
# original token), wrapped in literal "{" / "}" characters instead of the
# field delimiters, while the "_GoBack" bookmark around the tokens stays
# exactly where it was.
#
# Word's Field/Range character-position model treats fldChar/instrText runs
# as zero-width, so naive Range(start,end) pokes land on the wrong run. The
# reliable way to apply this kind of structural rewrite is to replace the
# whole host paragraph's Range with freshly authored OOXML via InsertXML.

$d = $word.ActiveDocument

# Locate the field and the paragraph that hosts it.
$field = $d.Fields.Item(1)
$fieldPos = $field.Code.Start

$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if (($fieldPos -ge $candidate.Range.Start) -and ($fieldPos -lt $candidate.Range.End)) {
        $targetParagraph = $candidate
    }
}

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:r><w:t>{</w:t></w:r>' `
    + '<w:r><w:t>m</w:t></w:r>' `
    + '<w:r><w:t>:</w:t></w:r>' `
    + '<w:r><w:t>''</w:t></w:r>' `
    + '<w:r><w:t>doc.html</w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
    + '<w:bookmarkEnd w:id="0"/>' `
    + '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">}</w:t></w:r>' `
    + '</w:p>'

$targetParagraph.Range.InsertXML($newParagraphXml)
